# Internet.xlsx update: add an "icons" sheet (icons used on each button),
# tidy up the "quickref.me" cell's number format (which collapses the
# duplicate "FF000000 w/o charset" font down onto the existing Link-less
# Arial/FF000000 font), and refresh the remembered cell selection on every
# sheet.

$wb = $excel.ActiveWorkbook

# --- fieldnames: remembered selection moves to B2 -------------------------
$wsFieldnames = $wb.Worksheets.Item("fieldnames")
$wsFieldnames.Range("B2").Select()

# --- URL: remembered selection moves to B2, tidy the quickref.me cell -----
$wsUrl = $wb.Worksheets.Item("URL")
$wsUrl.Range("C6").NumberFormat = "General"
$wsUrl.Range("B2").Select()

# --- color: remembered selection moves to B3 -------------------------------
$wsColor = $wb.Worksheets.Item("color")
$wsColor.Range("B3").Select()

# --- comments: remembered selection moves to B3 ----------------------------
$wsComments = $wb.Worksheets.Item("comments")
$wsComments.Range("B3").Select()

# --- icons: brand new sheet at the end of the tab strip --------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsIcons = $wb.Worksheets.Add($null, $lastSheet)
$wsIcons.Name = "icons"

$wsIcons.Range("A2").Value = "globe.png"
$wsIcons.Range("B2").Value = "BCN_Logo3.png"
$wsIcons.Range("C2").Value = "globe.png"
$wsIcons.Range("D2").Value = "info-circle.png"

$wsIcons.Range("A3").Value = "globe.png"
$wsIcons.Range("B3").Value = "BCN_Logo3.png"
$wsIcons.Range("C3").Value = "globe.png"

$wsIcons.Range("B4").Value = "info-circle.png"
$wsIcons.Range("C4").Value = "info-circle_red.png"

$wsIcons.Range("B5").Value = "info-circle.png"
$wsIcons.Range("C5").Value = "info-circle.png"

$wsIcons.Range("B6").Value = "info-circle.png"
$wsIcons.Range("C6").Value = "info-circle.png"

$wsIcons.Range("B7").Value = "info-circle.png"
$wsIcons.Range("C7").Value = "info-circle.png"

$wsIcons.Range("C8").Value = "info-circle.png"
$wsIcons.Range("C9").Value = "info-circle.png"

# icons becomes the active sheet/tab, selection parked on C5
$wsIcons.Range("C5").Select()
